$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply a Text number format to the whole ISBN column so the 13-digit
# codes are stored/rendered as text (keeps leading structure, matches
# the "Text" cellXf the workbook picks up).
$ws.Columns("A:A").NumberFormat = "@"

# Drop the two existing ISBN rows (old 10/13-digit hyphenated values) -
# they get re-entered below together with three brand new ISBN rows.
$ws.Rows(2).Delete()
$ws.Rows(2).Delete()

# Make sure the (re)created data rows keep the Text format too.
$ws.Range("A2:A6").NumberFormat = "@"

# Re-populate with the 5 ISBN-13 values (no dashes).
$ws.Range("A2").Value = "9780590353427"
$ws.Range("A3").Value = "9781338216660"
$ws.Range("A4").Value = "9780006479888"
$ws.Range("A5").Value = "9780141199702"
$ws.Range("A6").Value = "9780201835953"

# Move the selection below the new data, matching the authored workbook.
$ws.Range("A7").Select() | Out-Null
